$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 128.94737
$ws.Range("I39").Value = 71.90909000000001
$ws.Range("J39").Value = 207.375
$ws.Range("K39").Value = 215.72727
$ws.Range("L39").Value = 622.125
$ws.Range("M39").Value = 80.27272999999997
$ws.Range("N39").Value = -1214.125
$ws.Range("H44").Value = 14750
$ws.Range("J44").Value = 14750
$ws.Range("L44").Value = 14750
$ws.Range("N44").Value = -15674
$ws.Range("H64").Value = 3503.3704
$ws.Range("I64").Value = 3454.65
$ws.Range("J64").Value = 3642.5715
$ws.Range("K64").Value = 3454.65
$ws.Range("L64").Value = 3642.5715
$ws.Range("M64").Value = -3206.65
$ws.Range("N64").Value = -4138.5715
$ws.Range("H67").Value = 3503.3704
$ws.Range("I67").Value = 3454.65
$ws.Range("J67").Value = 3642.5715
$ws.Range("K67").Value = 3454.65
$ws.Range("L67").Value = 3642.5715
$ws.Range("M67").Value = -2596.65
$ws.Range("N67").Value = -5358.5715
$ws.Range("H74").Value = 3741.125
$ws.Range("I74").Value = 3405.1667
$ws.Range("J74").Value = 4749
$ws.Range("K74").Value = 3405.1667
$ws.Range("L74").Value = 4749
$ws.Range("M74").Value = -2469.1667
$ws.Range("N74").Value = -6621
$ws.Range("H76").Value = 3763
$ws.Range("I76").Value = 3668.6667
$ws.Range("J76").Value = 3819.6
$ws.Range("K76").Value = 3668.6667
$ws.Range("L76").Value = 3819.6
$ws.Range("M76").Value = -3353.6667
$ws.Range("N76").Value = -4449.6
$ws.Range("H77").Value = 3741.125
$ws.Range("I77").Value = 3405.1667
$ws.Range("J77").Value = 4749
$ws.Range("K77").Value = 17025.8335
$ws.Range("L77").Value = 23745
$ws.Range("M77").Value = -12345.8335
$ws.Range("N77").Value = -33105
$ws.Range("H79").Value = 3763
$ws.Range("I79").Value = 3668.6667
$ws.Range("J79").Value = 3819.6
$ws.Range("K79").Value = 3668.6667
$ws.Range("L79").Value = 3819.6
$ws.Range("M79").Value = -2576.6667
$ws.Range("N79").Value = -6003.6
$ws.Range("H86").Value = 50005148
$ws.Range("I86").Value = 500000000
$ws.Range("J86").Value = 5721.1113
$ws.Range("K86").Value = 500000000
$ws.Range("L86").Value = 5721.1113
$ws.Range("M86").Value = -499998877
$ws.Range("N86").Value = -7967.1113
$ws.Range("H89").Value = 50005148
$ws.Range("I89").Value = 500000000
$ws.Range("J89").Value = 5721.1113
$ws.Range("K89").Value = 2500000000
$ws.Range("L89").Value = 28605.5565
$ws.Range("M89").Value = -2499994384
$ws.Range("N89").Value = -39837.5565
$ws.Range("H107").Value = 14869.286
$ws.Range("I107").Value = 20702
$ws.Range("J107").Value = 287.5
$ws.Range("K107").Value = 20702
$ws.Range("L107").Value = 287.5
$ws.Range("M107").Value = -18782
$ws.Range("N107").Value = -4127.5
$ws.Range("H112").Value = 1908.6207
$ws.Range("J112").Value = 2005.7693
$ws.Range("L112").Value = 6017.3079
$ws.Range("N112").Value = -8233.3079
$ws.Range("H132").Value = 3678610.8
$ws.Range("I132").Value = 1913.4354
$ws.Range("J132").Value = 41671150
$ws.Range("K132").Value = 5740.3062
$ws.Range("L132").Value = 125013450
$ws.Range("M132").Value = -3210.3062
$ws.Range("N132").Value = -125018510

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1959.9445
$ws.Range("I61").Value = 1885.3269
$ws.Range("K61").Value = 1885.3269
$ws.Range("M61").Value = -1673.3269
$ws.Range("H63").Value = 35716950
$ws.Range("I63").Value = 41669260
$ws.Range("K63").Value = 41669260
$ws.Range("M63").Value = -41668574
$ws.Range("H66").Value = 35716950
$ws.Range("I66").Value = 41669260
$ws.Range("K66").Value = 208346300
$ws.Range("M66").Value = -208342868
$ws.Range("H136").Value = 1959.9445
$ws.Range("I136").Value = 1885.3269
$ws.Range("K136").Value = 5655.9807
$ws.Range("M136").Value = -3105.9807

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4102.5947
$ws.Range("I105").Value = 2735.8
$ws.Range("J105").Value = 4608.815
$ws.Range("K105").Value = 2735.8
$ws.Range("L105").Value = 4608.815
$ws.Range("M105").Value = -988.8000000000002
$ws.Range("N105").Value = -8102.815

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4939.521
$ws.Range("I31").Value = 1079.0714
$ws.Range("J31").Value = 7453.3022
$ws.Range("K31").Value = 1079.0714
$ws.Range("L31").Value = 7453.3022
$ws.Range("M31").Value = -784.0714
$ws.Range("N31").Value = -8043.3022
$ws.Range("H34").Value = 4939.521
$ws.Range("I34").Value = 1079.0714
$ws.Range("J34").Value = 7453.3022
$ws.Range("K34").Value = 1079.0714
$ws.Range("L34").Value = 7453.3022
$ws.Range("M34").Value = -877.0714
$ws.Range("N34").Value = -7857.3022
$ws.Range("H51").Value = 9596.571
$ws.Range("J51").Value = 9596.571
$ws.Range("L51").Value = 9596.571
$ws.Range("N51").Value = -11068.571
$ws.Range("H61").Value = 9596.571
$ws.Range("J61").Value = 9596.571
$ws.Range("L61").Value = 9596.571
$ws.Range("N61").Value = -10292.571
$ws.Range("H62").Value = 2878.375
$ws.Range("I62").Value = 2569.7058
$ws.Range("J62").Value = 3628
$ws.Range("K62").Value = 2569.7058
$ws.Range("L62").Value = 3628
$ws.Range("M62").Value = -1945.7058
$ws.Range("N62").Value = -4876
$ws.Range("H65").Value = 2878.375
$ws.Range("I65").Value = 2569.7058
$ws.Range("J65").Value = 3628
$ws.Range("K65").Value = 12848.529
$ws.Range("L65").Value = 18140
$ws.Range("M65").Value = -9728.529
$ws.Range("N65").Value = -24380

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1926.9375
$ws.Range("I132").Value = 1904
$ws.Range("J132").Value = 1928.4667
$ws.Range("K132").Value = 17136
$ws.Range("L132").Value = 17356.2003
$ws.Range("M132").Value = -14606
$ws.Range("N132").Value = -22416.2003

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 775
$ws.Range("I41").Value = 775
$ws.Range("K41").Value = 775
$ws.Range("M41").Value = -420
$ws.Range("H70").Value = 9958.177
$ws.Range("I70").Value = 13081
$ws.Range("J70").Value = 4233
$ws.Range("K70").Value = 13081
$ws.Range("L70").Value = 4233
$ws.Range("M70").Value = -12811
$ws.Range("N70").Value = -4773
$ws.Range("H73").Value = 9958.177
$ws.Range("I73").Value = 13081
$ws.Range("J73").Value = 4233
$ws.Range("K73").Value = 13081
$ws.Range("L73").Value = 4233
$ws.Range("M73").Value = -12145
$ws.Range("N73").Value = -6105
$ws.Range("H80").Value = 2879.0527
$ws.Range("I80").Value = 2711.6667
$ws.Range("J80").Value = 3029.7
$ws.Range("K80").Value = 2711.6667
$ws.Range("L80").Value = 3029.7
$ws.Range("M80").Value = -1713.6667
$ws.Range("N80").Value = -5025.7
$ws.Range("H83").Value = 2879.0527
$ws.Range("I83").Value = 2711.6667
$ws.Range("J83").Value = 3029.7
$ws.Range("K83").Value = 13558.3335
$ws.Range("L83").Value = 15148.5
$ws.Range("M83").Value = -8566.333500000001
$ws.Range("N83").Value = -25132.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25014344
$ws.Range("I132").Value = 9659.532999999999
$ws.Range("J132").Value = 100028400
$ws.Range("K132").Value = 28978.599
$ws.Range("L132").Value = 300085200
$ws.Range("M132").Value = -26448.599
$ws.Range("N132").Value = -300090260
